$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Rows 17 and 18 swap coin identity (Avalanche now appears before WrappedBTC)
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"

# Updated Price (D) and Volume(1h) (E) values for every data row
$updates = @{
    2  = @{ D = "30.169.79";     E = "  -4.10%  " }
    3  = @{ D = "1.912.73";      E = "  -3.87%  " }
    4  = @{ D = "0.9998";        E = "  +0.10%  " }
    5  = @{ D = "245.89";        E = "  -3.14%  " }
    6  = @{ D = "0.7013";        E = "  -14.53%  " }
    7  = @{ D = "1.002";         E = "  +0.25%  " }
    8  = @{ D = "0.3223";        E = "  -5.66%  " }
    9  = @{ D = "25.82";         E = "  +0.15%  " }
    10 = @{ D = "0.06832";       E = "  -2.72%  " }
    11 = @{ D = "0.7859";        E = "  -7.27%  " }
    12 = @{ D = "0.07934";       E = "  -2.30%  " }
    13 = @{ D = "1.914.04";      E = "  -3.80%  " }
    14 = @{ D = "5.364";         E = "  -2.37%  " }
    15 = @{ D = "93.71";         E = "  -8.17%  " }
    16 = @{ D = "258.84";        E = "  -5.88%  " }
    17 = @{ D = "14.30";         E = "  +2.07%  " }
    18 = @{ D = "30.165.66";     E = "  -4.12%  " }
    19 = @{ D = "5.769";         E = "  +0.67%  " }
    20 = @{ D = "0.000007836";   E = "  -1.87%  " }
    21 = @{ D = "1.001";         E = "  +0.18%  " }
    22 = @{ D = "2.167.34";      E = "  -3.49%  " }
    23 = @{ D = "0.9985";        E = "  -0.10%  " }
    24 = @{ D = "6.786";         E = "  -2.19%  " }
    25 = @{ D = "9.525";         E = "  -1.58%  " }
    26 = @{ D = "158.64";        E = "  -4.08%  " }
    27 = @{ D = "18.71";         E = "  -5.57%  " }
    28 = @{ D = "0.1302";        E = "  -16.36%  " }
    29 = @{ D = "2.199";         E = "  -0.22%  " }
    30 = @{ D = "1.365";         E = "  +0.76%  " }
    31 = @{ D = "1.546";         E = "  -1.57%  " }
    32 = @{ D = "4.389";         E = "  -4.14%  " }
    33 = @{ D = "4.159";         E = "  -3.99%  " }
    34 = @{ D = "0.05010";       E = "  -3.70%  " }
    35 = @{ D = "1.180";         E = "  -3.30%  " }
    36 = @{ D = "0.7378";        E = "  -1.85%  " }
    37 = @{ D = "2.726";         E = "  -1.70%  " }
    38 = @{ D = "0.01912";       E = "  -4.91%  " }
    39 = @{ D = "2.789";         E = "  -4.93%  " }
    40 = @{ D = "78.93";         E = "  +0.26%  " }
    41 = @{ D = "6.479";         E = "  -2.60%  " }
    42 = @{ D = "0.4398";        E = "  -6.06%  " }
    43 = @{ D = "1.997";         E = "  -3.94%  " }
    44 = @{ E = "  +0.29%  " }
    45 = @{ D = "0.8287";        E = "  -3.38%  " }
    46 = @{ D = "101.76";        E = "  -5.09%  " }
    47 = @{ D = "9.567";         E = "  -4.66%  " }
    48 = @{ D = "7.176";         E = "  -4.48%  " }
    49 = @{ D = "35.86";         E = "  -1.85%  " }
    50 = @{ D = "0.05889";       E = "  -1.61%  " }
    51 = @{ D = "1.460";         E = "  +1.32%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        Set-TextValue $ws.Range("D$row") $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
